# NIT-9016458460.xlsx — "Elimina antiguos EC y agrega nuevos y modifica Antigua BD"
#
# 1) Bump "Valor Mora" total (E11) and "Cant. Periodos" (F13).
# 2) Add a new period block ("2509") to the worker table, replicating the
#    previous period's ("2508") 9 worker rows (same B/C/D/F/G values & styles),
#    inserted right after the existing table and before the signature block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header totals -----------------------------------------------------
$ws.Range("E11").Value = 2726320
$ws.Range("F13").Value = 10

# --- insert 9 new rows for period 2509, copying the 2508 block ---------
$ws.Rows("54:62").Insert()

$srcBlock = $ws.Range("B45:J53")
$dstBlock = $ws.Range("B54:J62")
$srcBlock.Copy($dstBlock)

# Update the "Periodo Mora" column for the newly inserted rows
$ws.Range("E54:E62").Value = "2509"
